$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44326, 29533.279999999999, 30684.29),
    @(44327, 29537.09, 30687.41),
    @(44328, 29540.89, 30690.54),
    @(44329, 29544.69, 30693.66),
    @(44330, 29548.5, 30696.78),
    @(44331, 29552.3, 30699.91),
    @(44332, 29556.11, 30703.03),
    @(44333, 29559.919999999998, 30706.16),
    @(44334, 29563.72, 30709.279999999999),
    @(44335, 29567.53, 30712.41),
    @(44336, 29571.34, 30715.53),
    @(44337, 29575.15, 30718.66),
    @(44338, 29578.959999999999, 30721.78),
    @(44339, 29582.77, 30724.91),
    @(44340, 29586.57, 30728.03),
    @(44341, 29590.39, 30731.16),
    @(44342, 29594.2, 30734.29),
    @(44343, 29598.01, 30737.42),
    @(44344, 29601.82, 30740.54),
    @(44345, 29605.63, 30743.67),
    @(44346, 29609.439999999999, 30746.799999999999),
    @(44347, 29613.26, 30749.93),
    @(44348, 29617.07, 30753.06),
    @(44349, 29620.880000000001, 30756.19),
    @(44350, 29624.7, 30759.32),
    @(44351, 29628.51, 30762.45),
    @(44352, 29632.33, 30765.58),
    @(44353, 29636.15, 30768.71),
    @(44354, 29639.96, 30771.84),
    @(44355, 29643.78, 30774.97),
    @(44356, 29647.599999999999, 30778.1)
)

$startRow = 864
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Copy formatting from row 863 down to new rows
$ws.Range("A863:C863").Copy()
$ws.Range("A864:C894").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths: merge B:C into one width
$ws.Range("B2:C3").WrapText = $true
$ws.Rows.Item(2).RowHeight = 51

$ws.Columns.Item(2).ColumnWidth = 11.42578125
$ws.Columns.Item(3).ColumnWidth = 11.42578125

$wb.Names.Item("UF_IVP_DIARIO").RefersToR1C1Local = "=UF_IVP_DIARIO!R1C1:R894C3"
